$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "26.107.21"
Set-TextCell "E2" "  -0.79%  "
Set-TextCell "D3" "1.653.38"
Set-TextCell "E3" "  -0.82%  "
Set-TextCell "E4" "  -0.46%  "
Set-TextCell "D5" "218.62"
Set-TextCell "E5" "  -0.85%  "
Set-TextCell "D6" "0.5254"
Set-TextCell "E6" "  -0.97%  "
Set-TextCell "D8" "0.2669"
Set-TextCell "E8" "  +0.90%  "
Set-TextCell "D9" "0.06361"
Set-TextCell "E9" "  +0.05%  "
Set-TextCell "D10" "20.59"
Set-TextCell "E10" "  -1.57%  "
Set-TextCell "D11" "0.07702"
Set-TextCell "E11" "  -1.84%  "
Set-TextCell "B12" "Polkadot"
Set-TextCell "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D12" "4.601"
Set-TextCell "E12" "  +1.62%  "
Set-TextCell "B13" "WrappedEther"
Set-TextCell "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D13" "1.619.81"
Set-TextCell "E13" "  -3.06%  "
Set-TextCell "D14" "1.880.21"
Set-TextCell "E14" "  -0.82%  "
Set-TextCell "D15" "0.5608"
Set-TextCell "E15" "  +0.04%  "
Set-TextCell "D16" "0.0₅8222"
Set-TextCell "E16" "  +1.13%  "
Set-TextCell "D17" "65.40"
Set-TextCell "E17" "  -0.56%  "
Set-TextCell "D18" "26.112.54"
Set-TextCell "E18" "  -0.76%  "
Set-TextCell "D19" "1.004"
Set-TextCell "E19" "  -0.41%  "
Set-TextCell "D20" "4.702"
Set-TextCell "E20" "  -0.31%  "
Set-TextCell "D21" "10.39"
Set-TextCell "E21" "  +1.20%  "
Set-TextCell "D22" "191.14"
Set-TextCell "E22" "  -3.80%  "
Set-TextCell "D23" "5.985"
Set-TextCell "E23" "  -1.16%  "
Set-TextCell "E24" "  -0.49%  "
Set-TextCell "D25" "146.07"
Set-TextCell "E25" "  -0.70%  "
Set-TextCell "D26" "0.1201"
Set-TextCell "E26" "  -0.95%  "
Set-TextCell "D27" "7.260"
Set-TextCell "E27" "  +0.48%  "
Set-TextCell "D28" "15.94"
Set-TextCell "E28" "  -1.54%  "
Set-TextCell "D29" "1.496"
Set-TextCell "E29" "  -1.13%  "
Set-TextCell "D30" "0.05645"
Set-TextCell "E30" "  -3.98%  "
Set-TextCell "E31" "  -1.07%  "
Set-TextCell "D32" "3.506"
Set-TextCell "E32" "  -0.80%  "
Set-TextCell "D33" "3.388"
Set-TextCell "E33" "  +2.05%  "
Set-TextCell "D34" "1.583"
Set-TextCell "E34" "  -1.21%  "
Set-TextCell "D35" "2.796"
Set-TextCell "E35" "  -1.22%  "
Set-TextCell "D36" "0.9477"
Set-TextCell "E36" "  -1.33%  "
Set-TextCell "D37" "2.407"
Set-TextCell "E37" "  -1.00%  "
Set-TextCell "D38" "0.5788"
Set-TextCell "E38" "  -0.22%  "
Set-TextCell "D39" "0.01595"
Set-TextCell "D40" "5.976"
Set-TextCell "E40" "  +0.26%  "
Set-TextCell "E41" "  -0.53%  "
Set-TextCell "D42" "0.8403"
Set-TextCell "E42" "  -2.02%  "
Set-TextCell "D43" "1.021.12"
Set-TextCell "E43" "  -5.10%  "
Set-TextCell "D44" "101.67"
Set-TextCell "E44" "  -1.03%  "
Set-TextCell "D45" "1.791.48"
Set-TextCell "E45" "  -0.81%  "
Set-TextCell "D46" "58.44"
Set-TextCell "E46" "  +0.00%  "
Set-TextCell "B47" "Frax"
Set-TextCell "C47" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell "D47" "1.004"
Set-TextCell "E47" "  -0.95%  "
Set-TextCell "B48" "Cronos"
Set-TextCell "C48" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D48" "0.05339"
Set-TextCell "E48" "  +3.58%  "
Set-TextCell "B49" "BabyDogeCoin"
Set-TextCell "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D49" "0.0₈103"
Set-TextCell "E49" "  -1.65%  "
Set-TextCell "D50" "8.041"
Set-TextCell "E50" "  -0.08%  "
Set-TextCell "D51" "0.4343"
Set-TextCell "E51" "  -1.56%  "
